$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: find the paragraph whose text starts with a given marker string.
# ---------------------------------------------------------------------------
function Get-ParagraphStartingWith($doc, [string]$marker) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $para = $doc.Paragraphs.Item($i)
        if ($para.Range.Text.StartsWith($marker)) {
            return $para
        }
    }
    return $null
}

# ---------------------------------------------------------------------------
# 1) "Once complete the script will create a worklist ..." paragraph gains a
#    trailing " (Separate from the Daily List)".
# ---------------------------------------------------------------------------
$pWorklist = Get-ParagraphStartingWith $d "Once complete the script will create a worklist"
$rWorklist = $pWorklist.Range
$rWorklist.End = $rWorklist.End - 1
$rWorklist.Collapse(0)
$rWorklist.InsertAfter(" (Separate from the Daily List)")

# ---------------------------------------------------------------------------
# 2) "An Email will be sent ..." paragraph gains a trailing
#    " – include a link to the worklist." (en dash, not hyphen).
# ---------------------------------------------------------------------------
$pEmail = Get-ParagraphStartingWith $d "An Email will be sent to the QI member"
$rEmail = $pEmail.Range
$rEmail.End = $rEmail.End - 1
$rEmail.Collapse(0)
$enDash = [string][char]0x2013
$rEmail.InsertAfter(" " + $enDash + " include a link to the worklist.")

# ---------------------------------------------------------------------------
# 3) A brand-new third-level bullet is added right after that same paragraph:
#    "There will be a threshold on the number of cases on the worklist. If
#    the threshold is passed, the script will Email the QI team to request
#    support."  It sits at the same list (numId=1) but at ilvl=2 (Word's
#    1-based ListLevelNumber 3), matching sibling bullets like the "NOTE: ..."
#    paragraph further down.
# ---------------------------------------------------------------------------
$pEmail = Get-ParagraphStartingWith $d "An Email will be sent to the QI member"
$pEmail.Range.InsertParagraphAfter()

$pThreshold = Get-ParagraphStartingWith $d "An Email will be sent to the QI member"
$pThreshold = $pThreshold.Next()
$pThreshold.Range.ListFormat.ListLevelNumber = 3
$pThreshold.Range.Text = "There will be a threshold on the number of cases on the worklist. If the threshold is passed, the script will Email the QI team to request support."
